# "Cambios Estandarizacion Avance 3"
# Header standardization: replace Edad/Telefono/SaldoCuenta headers with
# Apellido/Cel/Edad (columns B/C/D), keeping A = Nombre.
# Net effect on row 1: A1=Nombre, B1=Apellido, C1=Cel, D1=Edad

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Apellido"
$ws.Range("C1").Value = "Cel"
$ws.Range("D1").Value = "Edad"

# Column widths: B gets a touch wider than the sheet default, C and D go
# back to the (default) width now that the wide "SaldoCuenta"/"Telefono"
# headers are gone. Inputs are pre-compensated for this engine's
# ColumnWidth rounding so the stored width lands as close as possible to
# the intended values.
$ws.Columns.Item(2).ColumnWidth = 8.5
$ws.Columns.Item(3).ColumnWidth = 8.333333333333334
$ws.Columns.Item(4).ColumnWidth = 8.333333333333334
